# Feature: Actions ahora chequea precios de polygon y envia mensajes al bot de Telegram
#
# Restructure the workbook:
#   1. Move "Totales" (currently last) to the front.
#   2. Insert a new, empty "Acciones invertidas" sheet right after "Totales".
#   3. Update the top-10-companies list on "Top 10 empresas del momento".
#   4. Update the low-price-stocks list on "Top 10 acciones bajas".
#
# Final tab order: Totales, Acciones invertidas, Top 10 empresas del momento,
# Top 10 acciones bajas.

$wb = $excel.ActiveWorkbook

# --- 1. Move "Totales" to be the first sheet -------------------------------
$totales = $wb.Worksheets.Item("Totales")
$totales.Move($wb.Worksheets.Item(1))

# --- 2. Insert new empty sheet "Acciones invertidas" right after it --------
$empresas = $wb.Worksheets.Item("Top 10 empresas del momento")
$accionesInvertidas = $wb.Worksheets.Add($empresas)
$accionesInvertidas.Name = "Acciones invertidas"

# --- 3. Update "Top 10 empresas del momento" contents -----------------------
$ws1 = $wb.Worksheets.Item("Top 10 empresas del momento")
$ws1.Range("A6").Value = "Google"
$ws1.Range("A7").Value = "Alibaba"
$ws1.Range("A8").Value = "Facebook"
$ws1.Range("A9").Value = "Berkshire Hathaway"
$ws1.Range("A10").Value = "Johnson & Johnson"
$ws1.Range("A11").Value = "Diversify your portfolio"

# --- 4. Update "Top 10 acciones bajas" contents (two identical blocks) -----
$ws2 = $wb.Worksheets.Item("Top 10 acciones bajas")

$ws2.Range("A2").Value = "PennyMac Financial Services"
$ws2.Range("A3").Value = "Amkor Technology"
$ws2.Range("A4").Value = "Air Lease Corporation"
$ws2.Range("A5").Value = "Baozun"
$ws2.Range("A6").Value = "Qurate Retail"
$ws2.Range("A7").Value = "Guess"
$ws2.Range("A8").Value = "UniFirst"
$ws2.Range("A9").Value = "Commvault Systems"
$ws2.Range("A10").Value = "Prudental Financial"
$ws2.Range("A11").Value = "Understanding and using eToro's CopyTrader feature wisely"

$ws2.Range("A13").Value = "PennyMac Financial Services"
$ws2.Range("A14").Value = "Amkor Technology"
$ws2.Range("A15").Value = "Air Lease Corporation"
$ws2.Range("A16").Value = "Baozun"
$ws2.Range("A17").Value = "Qurate Retail"
$ws2.Range("A18").Value = "Guess"
$ws2.Range("A19").Value = "UniFirst"
$ws2.Range("A20").Value = "Commvault Systems"
$ws2.Range("A21").Value = "Prudental Financial"
$ws2.Range("A22").Value = "Understanding and using eToro's CopyTrader feature wisely"
